# Update countries & provincias Spain
# - Refresh the "Datos actualizados..." timestamp in A1
# - Uzbekistan overtakes Chequia (rows 68/69 swap with updated counts)
# - Refresh case counts for several other countries (India, El Salvador,
#   Tailandia, Taiwan) without changing their row order
# NOTE: Groenlandia / Islas Malvinas (rows 209/210) swap display order in
# the source workbook purely via shared-string table reshuffling with no
# underlying numeric change (their B:H values are identical), so no cell
# write is required there to reproduce the same effective content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 08:15"

# Row 6 - India
$ws.Range("B6").Value = 879487
$ws.Range("C6").Value = 21
$ws.Range("E6").Value = 301864
$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 23194

# Row 68 - Uzbekistan overtakes Chequia, now shown first
$ws.Range("A68").Value = "Uzbekistan"
$ws.Range("B68").Value = 13193
$ws.Range("C68").Value = 196
$ws.Range("D68").Value = 7852
$ws.Range("E68").Value = 5280
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 61

# Row 69 - Chequia drops to second place (its own totals are unchanged,
# just relocated from the old row 68)
$ws.Range("A69").Value = "Chequia"
$ws.Range("B69").Value = 13174
$ws.Range("D69").Value = 8247
$ws.Range("E69").Value = 4575
$ws.Range("H69").Value = 352

# Row 76 - El Salvador
$ws.Range("D76").Value = 5663
$ws.Range("E76").Value = 3744
$ws.Range("G76").Value = 7
$ws.Range("H76").Value = 267

# Row 103 - Tailandia
$ws.Range("B103").Value = 3220
$ws.Range("C103").Value = 3
$ws.Range("D103").Value = 3090
$ws.Range("E103").Value = 72

# Row 159 - Taiwan
$ws.Range("D159").Value = 440
$ws.Range("E159").Value = 4
